$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Map-list table (rows 121-136) -----------------------------------------
# Add a new header row above the existing "map #" list, giving the three
# columns used below (bytes used w/o reprogram, bytes with subs, bytes with
# reprogram) real headers.
$ws.Range("B121").Value = "normal"
$ws.Range("C121").Value = "subs"
$ws.Range("D121").Value = "reprogram"

# Row 122 (map 1) stays as-is ("ok").

# Row 123 (map 2): "156 used" text -> numeric 156, with a comment.
$ws.Range("B123").Value = 156
$ws.Range("F123").Value = "annoying as hell"

# Row 124 (map 3): "110 used" text -> numeric 110; the old reprogram note
# ("62 with reprogram") moves from column C to column D, plus a comment.
$ws.Range("B124").Value = 110
$ws.Range("C124").Clear()
$ws.Range("D124").Value = "62 with reprogram"
$ws.Range("F124").Value = "annoying as hell"

# Row 125 (map 4): new data.
$ws.Range("B125").Value = 88
$ws.Range("F125").Value = "annoying as hell"
$ws.Range("K125").Value = "used left switch"

# Row 126 (map 5): new data.
$ws.Range("B126").Value = 92
$ws.Range("F126").Value = "annoying as hell"
$ws.Range("K126").Value = "fixed bad switch"

# Row 127 (map 6): new data.
$ws.Range("B127").Value = 78
$ws.Range("F127").Value = "very easy, should be an earlier level"

# Row 128 (map 7): new data.
$ws.Range("B128").Value = 64
$ws.Range("F128").Value = "very easy, should be an earlier level"

# Row 129 (map 8): new data.
$ws.Range("B129").Value = 156
$ws.Range("F129").Value = "fairly easy - fun water map"

# Row 130 (map 9): new data.
$ws.Range("B130").Value = 306
$ws.Range("C130").Value = 164
$ws.Range("F130").Value = "fun - lots of jumping"
$ws.Range("H130").Value = "jump-move forward and jump-move forward 3 are EASILY mass repeated on this map, will test with subs"
$ws.Range("I130").Value = "sub1 = jump/move forward, sub2 = jump"

# Row 131 (map 10): new data.
$ws.Range("B131").Value = 154
$ws.Range("F131").Value = "easy - very straightforward"

# Row 132 (map 11): new data.
$ws.Range("B132").Value = 128
$ws.Range("F132").Value = "not very complex, fairly easy as well - lots of random extra stuff not related to finishing the map"

# Row 133 (map 12): no new data beyond the existing map number.

# Row 134 (map 13): new data.
$ws.Range("B134").Value = 182
$ws.Range("F134").Value = "could probably lose the reprogram square, additionally, not sure if intended, but the last two switches can be skipped entirely, may want to disable jump on this map, or make some reason for the switches"

# Row 135 (map 14): new data.
$ws.Range("B135").Value = 318
$ws.Range("F135").Value = "very linear, interesting figuring out what does what"
$ws.Range("L135").Value = "needs edge squares removed maybe"

# Row 136 (map 15): new data (was a bare map-number row, now the last row).
$ws.Range("F136").Value = "can't be beaten"

# Rows 137-139 (maps 16-18) are removed entirely - the map list now stops at
# map 15/row 136.
$ws.Range("A137:P139").EntireRow.Delete()

# Reflect the new last-used cell as the active selection, matching the
# author's final cursor position after the edit.
$ws.Range("B136").Select()
